$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F3").Value = -3
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = -3
$ws.Range("F7").Value = -1
$ws.Range("F8").Value = -2
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = -2
$ws.Range("F18").Value = -4
$ws.Range("F19").Value = -7
$ws.Range("F20").Value = -5
$ws.Range("F21").Value = 7
$ws.Range("F23").Value = -6
$ws.Range("F24").Value = -4
$ws.Range("F25").Value = 0
$ws.Range("F28").Value = -4
$ws.Range("F31").Value = 1
$ws.Range("F32").Value = -8
$ws.Range("F33").Value = -4
$ws.Range("F34").Value = -2
$ws.Range("F36").Value = 4
$ws.Range("F37").Value = 3
$ws.Range("F39").Value = -4
$ws.Range("F41").Value = -4
$ws.Range("F42").Value = -5
$ws.Range("F45").Value = -3
$ws.Range("F48").Value = 1
